$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17:D17").NumberFormat = "@"

$ws.Range("A17").Value = "2025-01-10"
$ws.Range("B17").Value = "08:58:49"
$ws.Range("C17").Value = "Friday"
$ws.Range("D17").Value = "01"
$ws.Range("E17").Value = 126340
$ws.Range("F17").Value = 143710
$ws.Range("G17").Value = 169463
$ws.Range("H17").Value = 159509
$ws.Range("I17").Value = -1
$ws.Range("J17").Value = 142703
$ws.Range("K17").Value = -1
$ws.Range("L17").Value = -1
$ws.Range("M17").Value = 192809
$ws.Range("N17").Value = 115336
$ws.Range("O17").Value = 45784
$ws.Range("P17").Value = 28473
$ws.Range("Q17").Value = 64834
$ws.Range("R17").Value = -1
$ws.Range("S17").Value = 47897
$ws.Range("T17").Value = -1

$ws.Range("A17:T17").Style = "Normal"
